# Revert "Drop in files from RMI script"
# Updates the "About" sheet's reference-year text (2019 -> 2018) and the
# underlying conversion factor, which ripples through the three OCCF-Dp*
# sheets via their existing formulas.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

# Text labels that referenced "2019 dollars" should read "2018 dollars".
$about.Range("A18").Value = "billion 2018 dollars"
$about.Range("A21").Value = "million 2018 dollars"
$about.Range("A24").Value = "2018 dollars"
$about.Range("B26").Value = "2018 dollars per 2012 dollar"
$about.Range("B29").Value = 'which in this case is "2012 dollars per 2018 dollar."'

# Updated conversion factor (2012 dollars per 2018 dollar).
$about.Range("A26").Value = 0.9143273584567535

# Restore the default selection (A1) on the About sheet.
$about.Activate()
$about.Range("A1").Select()
